$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "switches" sheet right after "potentiometers" and make it active.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "switches"

# Row 1 headers - same header strings/styles as sheet1's row 1.
$headers = @{
    "A1" = "TPN"
    "B1" = "Description"
    "C1" = "Value"
    "D1" = "Tolerance"
    "E1" = "Taper"
    "F1" = "Power (Watts)"
    "G1" = "Temperature Coefficient"
    "H1" = "Size / Dimension"
    "I1" = "Mounting Type"
    "J1" = "Library Ref"
    "K1" = "Footprint Ref"
    "L1" = "Manufacturer 1"
    "M1" = "Manufacturer 1 PN"
    "N1" = "Supplier 1"
    "O1" = "Supplier 1 PN"
    "P1" = "Supplier 1 Link"
    "Q1" = "Manufacturer 2"
    "R1" = "Manufacturer 2 PN"
    "S1" = "Supplier 2"
    "T1" = "Supplier 2 PN"
}

foreach ($addr in $headers.Keys) {
    $cell = $ws2.Range($addr)
    $cell.Value = $headers[$addr]
    $cell.Font.Bold = $true
}

# Columns C (Value) and I (Mounting Type) are text-formatted, like on sheet1.
$ws2.Range("C1").NumberFormat = "@"
$ws2.Range("I1").NumberFormat = "@"

# Row 2: generated TPN via formula.
$ws2.Range("A2").Formula = '="SW-"&TEXT(ROW()-1,"0000")'

# Row 3: a second TPN entered as a literal value.
$ws2.Range("A3").Value = "SW-0002"

# Selection / active cell bookkeeping to match the saved UI state.
$ws1.Range("A2").Select() | Out-Null
$ws2.Range("A3").Select() | Out-Null
$ws2.Activate() | Out-Null
